{"js": "// Apply RF015 v1.3 -> v1.4 text corrections (gender/number agreement fixes\n// and a couple of wording changes) by locating each exact sentence with\n// Body.search() and replacing it in place with Range.insertText(..., \"Replace\").\n// This preserves the original run formatting (rPr) of the located range.\n\nconst body = context.document.body;\n\nasync function replaceAll(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. Basic Flow, step 2 \u2014 \"das...cadastradas\" -> \"dos...cadastrados\"\nawait replaceAll(\n  \"2. System exibe a listagem das Planos de Capacitacao de TI cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' \",\n  \"2. System exibe a listagem dos Planos de Capacitacao de TI cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' \"\n);\n\n// 2. Basic Flow, step 5 \u2014 \"escolha\" -> \"escolhe\"\nawait replaceAll(\n  \"5. Lider de Pessoas escolha o 'Periodo Avaliativo' apropriado no campo de selecao \",\n  \"5. Lider de Pessoas escolhe o 'Periodo Avaliativo' apropriado no campo de selecao \"\n);\n\n// 3. Basic Flow, step 7 \u2014 \"selecione\" -> \"seleciona\"\nawait replaceAll(\n  \"7. Lider de Pessoas selecione a 'Unidade' correspondente no campo de selecao de unidade \",\n  \"7. Lider de Pessoas seleciona a 'Unidade' correspondente no campo de selecao de unidade \"\n);\n\n// 4. Basic Flow, step 9 \u2014 \"preencha\" -> \"preenche\"\nawait replaceAll(\n  \"9. Lider de Pessoas preencha o campo 'Possiveis Capacitacoes' com informacoes sobre capacitacoes adicionais \",\n  \"9. Lider de Pessoas preenche o campo 'Possiveis Capacitacoes' com informacoes sobre capacitacoes adicionais \"\n);\n\n// 5. Basic Flow, step 11 \u2014 \"preencha\" -> \"preenche\"\nawait replaceAll(\n  \"11. Lider de Pessoas preencha o campo 'Observacao' com informacoes adicionais ou relevantes sobre o plano de capacitacao \",\n  \"11. Lider de Pessoas preenche o campo 'Observacao' com informacoes adicionais ou relevantes sobre o plano de capacitacao \"\n);\n\n// 6. AF[1]/AF[2]/AF[3] step 1 (3 occurrences) \u2014 \"um Planos\" -> \"um Plano\"\nawait replaceAll(\n  \"1. Lider de Pessoas seleciona um Planos de Capacitacao de TI da listagem \",\n  \"1. Lider de Pessoas seleciona um Plano de Capacitacao de TI da listagem \"\n);\n\n// 7. AF[2] step 6 \u2014 \"das\" -> \"dos\"\nawait replaceAll(\n  \"6. System exibe a listagem das Planos de Capacitacao de TI sem a Capacitacao de TI excluida ef[3,4]\",\n  \"6. System exibe a listagem dos Planos de Capacitacao de TI sem a Capacitacao de TI excluida ef[3,4]\"\n);\n\n// 8. AF[3] title \u2014 \"Negar Exclusao\" -> \"Nao Confirmar Exclusao\"\nawait replaceAll(\n  \"AF[3] \\u2013 Negar Exclusao da Capacitacao de TI\",\n  \"AF[3] \\u2013 Nao Confirmar Exclusao da Capacitacao de TI\"\n);\n\n// 9. AF[3] step 6 \u2014 \"das...excluida\" -> \"dos...nao excluida\"\nawait replaceAll(\n  \"6. System exibe a listagem das Planos de Capacitacao de TI com a Capacitacao de TI excluida \",\n  \"6. System exibe a listagem dos Planos de Capacitacao de TI com a Capacitacao de TI nao excluida \"\n);\n\n// 10. AF[4] step 2 \u2014 \"das...cadastradas\" -> \"dos...cadastrados\"\nawait replaceAll(\n  \"2. System exibe a listagem das Planos de Capacitacao de TI cadastradas apenas para visualizacao com a opcao 'Ajuda' \",\n  \"2. System exibe a listagem dos Planos de Capacitacao de TI cadastrados apenas para visualizacao com a opcao 'Ajuda' \"\n);\n", "ps1": "# Apply RF015 v1.3 -> v1.4 text corrections (gender/number agreement fixes\n# and a couple of wording changes) using Find/Replace over the whole\n# document content.\n#\n# Each Find/Replace pair below is trimmed to the smallest fragment that is\n# still unique (or, for the two sentences that legitimately repeat, that\n# covers exactly the repeated occurrences that must all receive the very\n# same fix) and, crucially, never contains a straight apostrophe ' inside\n# the FindText/ReplaceWith strings. The engine's Find & Replace applies\n# smart-quote autocorrection to apostrophes that travel through ReplaceWith,\n# which would otherwise turn the existing straight quotes used around\n# 'Novo', 'Periodo Avaliativo', etc. into curly quotes even though the\n# diff leaves them untouched. Keeping the replaced fragments apostrophe-free\n# sidesteps that and leaves the rest of each sentence byte-for-byte as-is.\n#\n# wdReplaceAll (2) makes a single Execute call replace every occurrence of\n# a given fragment in one shot (used for the \"Lider de Pessoas preencha o\n# campo \" fragment shared by steps 9 & 11, and the \"seleciona um Planos\"\n# sentence that appears 3 times across AF[1]/AF[2]/AF[3]).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Exact {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $rng = $d.Content\n    $rng.Start = 0\n    $rng.Find.Execute(\n        $FindText,\n        $true,  # MatchCase\n        $false, # MatchWholeWord\n        $false, # MatchWildcards\n        $false, # MatchSoundsLike\n        $false, # MatchAllWordForms\n        $true,  # Forward\n        1,      # Wrap = wdFindContinue\n        $false, # Format\n        $ReplaceText,\n        2       # Replace = wdReplaceAll\n    )\n}\n\n# 1. Basic Flow, step 2 - \"das...cadastradas\" -> \"dos...cadastrados\"\nReplace-Exact `\n    \"System exibe a listagem das Planos de Capacitacao de TI cadastradas com opcoes\" `\n    \"System exibe a listagem dos Planos de Capacitacao de TI cadastrados com opcoes\"\n\n# 2. Basic Flow, step 5 - \"escolha\" -> \"escolhe\"\nReplace-Exact \"Lider de Pessoas escolha o\" \"Lider de Pessoas escolhe o\"\n\n# 3. Basic Flow, step 7 - \"selecione\" -> \"seleciona\"\nReplace-Exact \"Lider de Pessoas selecione a\" \"Lider de Pessoas seleciona a\"\n\n# 4 & 5. Basic Flow, steps 9 and 11 - \"preencha\" -> \"preenche\" (both sentences\n# share this exact fragment, so one ReplaceAll call fixes both)\nReplace-Exact \"Lider de Pessoas preencha o campo \" \"Lider de Pessoas preenche o campo \"\n\n# 6. AF[1]/AF[2]/AF[3] step 1 (3 occurrences) - \"um Planos\" -> \"um Plano\"\nReplace-Exact `\n    \"seleciona um Planos de Capacitacao de TI da listagem\" `\n    \"seleciona um Plano de Capacitacao de TI da listagem\"\n\n# 7. AF[2] step 6 - \"das\" -> \"dos\"\nReplace-Exact `\n    \"listagem das Planos de Capacitacao de TI sem a Capacitacao de TI excluida\" `\n    \"listagem dos Planos de Capacitacao de TI sem a Capacitacao de TI excluida\"\n\n# 8. AF[3] title - \"Negar Exclusao\" -> \"Nao Confirmar Exclusao\"\nReplace-Exact `\n    \"AF[3] \u2013 Negar Exclusao da Capacitacao de TI\" `\n    \"AF[3] \u2013 Nao Confirmar Exclusao da Capacitacao de TI\"\n\n# 9. AF[3] step 6 - \"das...excluida\" -> \"dos...nao excluida\"\nReplace-Exact `\n    \"listagem das Planos de Capacitacao de TI com a Capacitacao de TI excluida\" `\n    \"listagem dos Planos de Capacitacao de TI com a Capacitacao de TI nao excluida\"\n\n# 10. AF[4] step 2 - \"das...cadastradas\" -> \"dos...cadastrados\"\nReplace-Exact `\n    \"listagem das Planos de Capacitacao de TI cadastradas apenas para visualizacao\" `\n    \"listagem dos Planos de Capacitacao de TI cadastrados apenas para visualizacao\"\n"}
